$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GmailID")
$ws.Range("B2").Value = "dlkzgzdeizmgpqje"
